$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.689.18'
$ws.Range("E2").Value = '  -2.50%  '
$ws.Range("D3").Value = '2.916.90'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '500.31'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.13'
$ws.Range("E6").Value = '  -4.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.13'
$ws.Range("E9").Value = '  -4.73%  '
$ws.Range("E10").Value = '  -4.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.350'
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").Value = '3.413.48'
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("E13").Value = '  -3.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.59'
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000158'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '55.678.06'
$ws.Range("E16").Value = '  -2.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.93'
$ws.Range("E17").Value = '  -4.15%  '
$ws.Range("D18").Value = '2.915.83'
$ws.Range("E18").Value = '  -2.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.58'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.68'
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '313.59'
$ws.Range("E21").Value = '  -4.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.70'
$ws.Range("E24").Value = '  -2.49%  '
$ws.Range("D25").Value = '3.031.52'
$ws.Range("E25").Value = '  -3.49%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -4.74%  '
$ws.Range("D28").Value = '0.0₃0833'
$ws.Range("E28").Value = '  -8.28%  '
$ws.Range("E29").Value = '  -6.27%  '
$ws.Range("E30").Value = '  -7.62%  '
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.84'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.06'
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.57'
$ws.Range("E36").Value = '  -4.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.93'
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("E38").Value = '  -6.36%  '
$ws.Range("E39").Value = '  -4.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.36'
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.635'
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").Value = '2.110.09'
$ws.Range("E44").Value = '  -7.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.94'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.32'
$ws.Range("E46").Value = '  -5.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.916'
$ws.Range("E47").Value = '  -6.10%  '
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.55'
$ws.Range("E49").Value = '  -3.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0836'
$ws.Range("E50").Value = '  -6.00%  '
$ws.Range("E51").Value = '  -8.84%  '
